# Adds server-side support for nitro boost detection to the "netcalls"
# table: three new KitEx.NetCalls rows (IDs 1138-1140) inserted right
# after "SendKitsAccess" (row 73), pushing every later row down by three.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows starting at row 74 (old rows 74.. shift to 77..).
$ws.Range("A74:A76").EntireRow.Insert()

# Fill the new rows in the same cell-write order the original author used,
# so new shared-string entries land in the same order/index as the source.
$ws.Cells.Item(74, 5).Value = "ulong[] players"
$ws.Cells.Item(76, 1).Value = "RequestIsNitroBoosting"
$ws.Cells.Item(74, 1).Value = "SendNitroBoostingUpdated"
$ws.Cells.Item(75, 1).Value = "RespondIsNitroBoosting"
$ws.Cells.Item(75, 5).Value = "byte[] response (0=not, 1=is, 2=unclear)"
$ws.Cells.Item(76, 5).Value = "ulong[] player, byte[] code (0=not, 1=is, 2=unclear)"

$ws.Cells.Item(74, 2).Value = 1138
$ws.Cells.Item(74, 3).Value = "KitEx.NetCalls"
$ws.Cells.Item(74, 4).Value = "FROM_CLIENT"

$ws.Cells.Item(75, 2).Value = 1139
$ws.Cells.Item(75, 3).Value = "KitEx.NetCalls"
$ws.Cells.Item(75, 4).Value = "FROM_SERVER"

$ws.Cells.Item(76, 2).Value = 1140
$ws.Cells.Item(76, 3).Value = "KitEx.NetCalls"
$ws.Cells.Item(76, 4).Value = "FROM_SERVER"

# Grow the "Table2" ListObject (and its AutoFilter range) to cover the three
# newly-added rows, so the table keeps wrapping the full A1:E114 data block.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E114"))

# Restore the on-screen selection to where the edit actually happened.
$ws.Range("C76").Select()
